# "nha sach chuan + mau menu"
# Insert a new "Giấy" sheet right after "Dụng cụ học sinh", populated like
# that first sheet but with a plain running number (1-10) in column B
# instead of the THSxx product-code text, then restore/update the various
# sheet-view selections (active tab moves to the new sheet, the old first
# sheet keeps only a single-cell selection, and the last sheet gets its
# scroll position + selection updated).

$wb = $excel.ActiveWorkbook

# Sheet that everything is modelled after / inserted next to.
$wsTools = $wb.Worksheets.Item(1)

# --- Create the new "Giấy" sheet right after "Dụng cụ học sinh" ----------
$wsPaper = $wb.Worksheets.Add($null, $wsTools)
$wsPaper.Name = "Giấy"

# Copy the header + data block (values, number formats, styles) straight
# from the "Dụng cụ học sinh" sheet so fonts/fills/borders/number formats
# all match, without bringing over its custom column widths.
$wsTools.Range("A1:E11").Copy($wsPaper.Range("A1"))

# Column B on the new sheet is a plain sequential number (1-10), not the
# THSxx product-code text used on the "Dụng cụ học sinh" sheet.
for ($i = 2; $i -le 11; $i++) {
    $wsPaper.Cells.Item($i, 2).Value2 = $i - 1
}

# --- Fix up view/selection state -----------------------------------------
# "Dụng cụ học sinh" no longer is the active tab, and its remembered
# selection shrinks down to the single cell D2.
$wsTools.Range("D2").Select()

# New sheet becomes the active tab, with its own remembered selection.
$wsPaper.Range("F22").Select()

# "Sách tham khảo" (now pushed one slot later by the insertion) gets a new
# scroll position/selection recorded in its sheet view.
$wsRef = $wb.Worksheets.Item("Sách tham khảo")
$wsRef.Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$wsRef.Range("K18").Select()

# Leave the new "Giấy" sheet focused/active, matching activeTab="1".
$wsPaper.Select()
